# Saldo.xlsx update — refresh the "Export" balance snapshot:
#   - remove stale/duplicated account rows
#   - insert the newly reported account balances in their sorted (descending
#     Saldo) position, matching the refreshed export

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Remove-AccountRow([string]$account) {
    $found = $ws.Cells.Find($account)
    $r = $found.Row()
    $ws.Rows.Item($r).Delete()
}

function Set-AccountRow([int]$r, [string]$account, [string]$name, [double]$saldo) {
    # Force column A to be stored as text so the zero-padded account number
    # survives (otherwise Excel's auto-detection would coerce "004207955"
    # into the number 4207955), then drop the leftover "@" number-format so
    # the cell's style matches its plain, unstyled neighbours.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $account
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $saldo
}

function Insert-AccountRowBefore([string]$beforeAccount, [string]$account, [string]$name, [double]$saldo) {
    $anchor = $ws.Cells.Find($beforeAccount)
    $r = $anchor.Row()
    $ws.Rows.Item($r).Insert()
    Set-AccountRow $r $account $name $saldo
}

function Insert-AccountRowAfter([string]$afterAccount, [string]$account, [string]$name, [double]$saldo) {
    $anchor = $ws.Cells.Find($afterAccount)
    $r = $anchor.Row() + 1
    $ws.Rows.Item($r).Insert()
    Set-AccountRow $r $account $name $saldo
}

# --- Removals (done first, while every account number is still unique) ---
Remove-AccountRow "004879567"   # SANDRA     95201.98
Remove-AccountRow "004567324"   # FRANCISCO  64986.65
Remove-AccountRow "004268684"   # PATRICIA   15.41 (stale duplicate balance)

# --- Additions, inserted to preserve the descending-Saldo sort order ---
Insert-AccountRowBefore "001731007" "004207955" "SILVANIA" 45.03
Insert-AccountRowBefore "005000656" "004001621" "DANIELA"  37.58
Insert-AccountRowBefore "005000656" "004268684" "PATRICIA" 37.29
Insert-AccountRowBefore "004340984" "004211922" "CARLOS"   34.71
Insert-AccountRowBefore "004643153" "004377415" "ANGELA"   26.37
Insert-AccountRowBefore "004381194" "004756968" "DANIELY"  18.08
Insert-AccountRowBefore "004453302" "004212581" "MARIA"    0.59
Insert-AccountRowAfter  "004371857" "004332207" "IRACY"    0.16
